$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)

# 1) "Text Placeholder 5": "for-цикъл" -> "foreach-цикъл"
$shpTitle = $s.Shapes.Item(1)
$trTitle = $shpTitle.TextFrame.TextRange
$subTitle = $trTitle.Characters(6, 9)
$subTitle.Text = "foreach-цикъл"

# 2) "Rectangle 4" code box: update the for-loop snippet to a foreach snippet.
$shpCode = $s.Shapes.Item(4)
$trCode = $shpCode.TextFrame.TextRange

# Replace the later substring first so the earlier substring's offsets stay valid.
$subLine3 = $trCode.Characters(77, 28)
$subLine3.Text = "Console.WriteLine(num); }"

$subLine2 = $trCode.Characters(36, 38)
$subLine2.Text = "foreach (int num in arr) {"

# The shape auto-fits its height to the (now shorter) text; pin it to the
# exact value PowerPoint's text layout produced for the edited snippet.
$shpCode.Height = 128.62185551181102
